$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record needs to be inserted as the new row 227 of the
# "Coliflor" price history. Every existing record from row 227 down to the
# previous last row (290) shifts down by one row (to 228..291) and this
# fresh record takes the now-empty row 227.
$ws.Rows.Item(227).Insert()

$row = 227
$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"
$ws.Cells.Item($row, 4).Value = 44642
$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = 100112008
$ws.Cells.Item($row, 7).Value = "Coliflor"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Segunda"
$ws.Cells.Item($row, 10).Value = 500
$ws.Cells.Item($row, 11).Value = 1700
$ws.Cells.Item($row, 12).Value = 1700
$ws.Cells.Item($row, 13).Value = 1700
$ws.Cells.Item($row, 14).Value = "$/unidad"
$ws.Cells.Item($row, 15).Value = "Región Metropolitana"
$ws.Cells.Item($row, 16).Value = 1700
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"
